$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number would be auto-converted
# from text to a numeric type by Excel on assignment (exactly like typing into
# a General-formatted cell). The source data keeps these as literal text
# (e.g. trailing zeros / scientific-looking tiny decimals must render verbatim),
# so we mark those specific cells as Text before writing the new value.
$textForceCells = @(
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.809.09"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "1.904.33"
$ws.Range("E3").Value = "  -0.26%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "312.81"
$ws.Range("E5").Value = "  -1.05%  "

$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").Value = "0.5037"
$ws.Range("E7").Value = "  +4.51%  "

$ws.Range("D8").Value = "0.3804"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "0.07269"
$ws.Range("E9").Value = "  -1.20%  "

$ws.Range("D10").Value = "0.9040"
$ws.Range("E10").Value = "  -3.14%  "

$ws.Range("D11").Value = "20.82"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "0.07664"
$ws.Range("E12").Value = "  -1.25%  "

$ws.Range("D13").Value = "1.865.05"
$ws.Range("E13").Value = "  -2.16%  "

$ws.Range("D14").Value = "5.470"
$ws.Range("E14").Value = "  -0.55%  "

$ws.Range("D15").Value = "91.57"
$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").Value = "0.000008699"
$ws.Range("E17").Value = "  -1.55%  "

$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("D19").Value = "27.836.70"
$ws.Range("E19").Value = "  -0.82%  "

$ws.Range("D20").Value = "14.54"
$ws.Range("E20").Value = "  -1.70%  "

$ws.Range("D21").Value = "5.160"
$ws.Range("E21").Value = "  -0.38%  "

$ws.Range("D22").Value = "10.80"
$ws.Range("E22").Value = "  -0.94%  "

$ws.Range("D23").Value = "6.561"
$ws.Range("E23").Value = "  -1.04%  "

$ws.Range("D24").Value = "153.79"
$ws.Range("E24").Value = "  -1.44%  "

$ws.Range("D25").Value = "1.876"
$ws.Range("E25").Value = "  -2.05%  "

$ws.Range("D26").Value = "2.215"
$ws.Range("E26").Value = "  +3.76%  "

$ws.Range("D27").Value = "18.34"
$ws.Range("E27").Value = "  -0.97%  "

$ws.Range("D28").Value = "115.18"
$ws.Range("E28").Value = "  -1.48%  "

$ws.Range("D29").Value = "4.897"
$ws.Range("E29").Value = "  -1.46%  "

$ws.Range("D30").Value = "0.09008"
$ws.Range("E30").Value = "  +0.66%  "

$ws.Range("D31").Value = "3.212"
$ws.Range("E31").Value = "  -2.99%  "

$ws.Range("D32").Value = "1.216"
$ws.Range("E32").Value = "  -3.09%  "

$ws.Range("D33").Value = "4.653"
$ws.Range("E33").Value = "  -0.57%  "

$ws.Range("D34").Value = "0.7607"

$ws.Range("D35").Value = "0.02061"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "2.496"
$ws.Range("E36").Value = "  -5.50%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "1.094"
$ws.Range("E37").Value = "  -1.49%  "

$ws.Range("D38").Value = "0.5506"
$ws.Range("E38").Value = "  +0.20%  "

$ws.Range("E39").Value = "  +0.72%  "

$ws.Range("D40").Value = "0.05234"
$ws.Range("E40").Value = "  -1.44%  "

$ws.Range("D41").Value = "6.849"
$ws.Range("E41").Value = "  -2.63%  "

$ws.Range("D42").Value = "8.435"
$ws.Range("E42").Value = "  -0.92%  "

$ws.Range("D43").Value = "0.1508"
$ws.Range("E43").Value = "  -1.55%  "

$ws.Range("D44").Value = "110.72"
$ws.Range("E44").Value = "  +2.47%  "

$ws.Range("D45").Value = "10.56"
$ws.Range("E45").Value = "  -1.42%  "

$ws.Range("D46").Value = "0.4791"
$ws.Range("E46").Value = "  -0.94%  "

$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("D48").Value = "1.621"
$ws.Range("E48").Value = "  -1.81%  "

$ws.Range("D49").Value = "67.13"
$ws.Range("E49").Value = "  -1.32%  "

$ws.Range("D50").Value = "0.06062"
$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("D51").Value = "0.9011"
$ws.Range("E51").Value = "  +0.07%  "
